$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D10").Value = 22
